$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Test" column (F) header + rows 1-4, plus a new "Comment" value in E3
# (replacing the old "String" placeholder). Author typed F1, F2, F3, F4 then
# E3, in that order (matches the resulting shared-string allocation order).
$ws.Range("F1").Value = "Test"
$ws.Range("F1").Font.Name = "等线"
$ws.Range("F1").HorizontalAlignment = -4131

$ws.Range("F2").Value = "Function"
$ws.Range("F3").Value = "Lua"
$ws.Range("F4").Value = "function() print(12) end"
$ws.Range("E3").Value = "Comment"

# Re-use F1's freshly created style (font + left alignment) for the rest of
# the new cells instead of re-deriving the font from scratch on each one.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("F2").PasteSpecial(-4122) | Out-Null
$ws.Range("F3").PasteSpecial(-4122) | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("E9").Select() | Out-Null
